$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '81.217.95'
$ws.Range("E2").Value = '  +5.03%  '

# Row 3
$ws.Range("D3").Value = '3.190.61'
$ws.Range("E3").Value = '  +1.62%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").Value = '''209.61'
$ws.Range("E5").Value = '  +3.54%  '

# Row 6
$ws.Range("D6").Value = '''634.65'
$ws.Range("E6").Value = '  +0.96%  '

# Row 7
$ws.Range("D7").Value = '''0.289'
$ws.Range("E7").Value = '  +27.93%  '

# Row 8
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '''0.594'
$ws.Range("E9").Value = '  +4.24%  '

# Row 10
$ws.Range("D10").Value = '3.192.98'
$ws.Range("E10").Value = '  +1.78%  '

# Row 11
$ws.Range("D11").Value = '''0.592'
$ws.Range("E11").Value = '  +11.93%  '

# Row 12
$ws.Range("D12").Value = '''0.0000265'
$ws.Range("E12").Value = '  +18.41%  '

# Row 13
$ws.Range("E13").Value = '  +2.31%  '

# Row 14
$ws.Range("D14").Value = '''5.42'
$ws.Range("E14").Value = '  +0.63%  '

# Row 15
$ws.Range("D15").Value = '3.787.98'
$ws.Range("E15").Value = '  +2.21%  '

# Row 16
$ws.Range("D16").Value = '''32.26'
$ws.Range("E16").Value = '  +6.12%  '

# Row 17
$ws.Range("D17").Value = '81.190.21'
$ws.Range("E17").Value = '  +5.28%  '

# Row 18
$ws.Range("D18").Value = '3.196.59'
$ws.Range("E18").Value = '  +2.26%  '

# Row 19
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").Value = '''3.26'
$ws.Range("E19").Value = '  +14.48%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''14.44'
$ws.Range("E20").Value = '  +3.00%  '

# Row 21
$ws.Range("D21").Value = '''9.28'
$ws.Range("E21").Value = '  +1.02%  '

# Row 22
$ws.Range("D22").Value = '''442.35'
$ws.Range("E22").Value = '  +2.59%  '

# Row 23
$ws.Range("D23").Value = '''5.25'
$ws.Range("E23").Value = '  +11.54%  '

# Row 24
$ws.Range("D24").Value = '''7.09'
$ws.Range("E24").Value = '  +5.48%  '

# Row 25
$ws.Range("E25").Value = '  +10.32%  '

# Row 26
$ws.Range("D26").Value = '''11.31'
$ws.Range("E26").Value = '  +7.38%  '

# Row 27
$ws.Range("D27").Value = '3.364.85'
$ws.Range("E27").Value = '  +2.21%  '

# Row 28
$ws.Range("D28").Value = '''77.31'
$ws.Range("E28").Value = '  +3.11%  '

# Row 29
$ws.Range("E29").Value = '  +13.52%  '

# Row 30
$ws.Range("D30").Value = '''0.993'
$ws.Range("E30").Value = '  -0.55%  '

# Row 31
$ws.Range("D31").Value = '''9.23'
$ws.Range("E31").Value = '  +5.70%  '

# Row 32
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.49%  '

# Row 33
$ws.Range("D33").Value = '''571.29'
$ws.Range("E33").Value = '  +9.16%  '

# Row 34
$ws.Range("D34").Value = '''1.52'
$ws.Range("E34").Value = '  +3.28%  '

# Row 35
$ws.Range("D35").Value = '''0.155'
$ws.Range("E35").Value = '  +14.85%  '

# Row 36
$ws.Range("D36").Value = '''2.04'
$ws.Range("E36").Value = '  +4.60%  '

# Row 37
$ws.Range("D37").Value = '''0.133'
$ws.Range("E37").Value = '  +24.95%  '

# Row 38
$ws.Range("D38").Value = '''23.15'
$ws.Range("E38").Value = '  +4.95%  '

# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  +0.10%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '''0.416'
$ws.Range("E40").Value = '  +6.69%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '''3.14'
$ws.Range("E41").Value = '  +24.73%  '

# Row 42
$ws.Range("D42").Value = '''5.99'
$ws.Range("E42").Value = '  +12.10%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''2.05'
$ws.Range("E43").Value = '  +17.00%  '

# Row 44
$ws.Range("D44").Value = '''20.80'
$ws.Range("E44").Value = '  +3.72%  '

# Row 45
$ws.Range("D45").Value = '''160.13'
$ws.Range("E45").Value = '  -2.20%  '

# Row 46
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''190.66'
$ws.Range("E47").Value = '  -2.10%  '

# Row 48
$ws.Range("D48").Value = '''1.35'
$ws.Range("E48").Value = '  +5.66%  '

# Row 49
$ws.Range("D49").Value = '''0.783'
$ws.Range("E49").Value = '  -1.39%  '

# Row 50
$ws.Range("D50").Value = '''43.64'
$ws.Range("E50").Value = '  +2.43%  '

# Row 51
$ws.Range("D51").Value = '''4.33'
$ws.Range("E51").Value = '  +5.92%  '
